$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHP5YWAN5%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T200014Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCu5FJXG3S5aKeF7FyN2EBzNWsS9JBYsFDsYKJGuEbZEQIgAhZCaqLzvYqjyzdSgwHh1ncWapjmeSRzcn9wOTShDXgq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDPKRBg3%2B7gik4%2FdoJCq6A60R3mr2sOyO17YUbq%2Bv4aDc365v9rrCE13OkzA4%2BfkMac6Rf4uuQROrUAtfb%2Bjp%2BN1bGN7HKqLD9xOmXEdNljVWjute3j05HX1X%2BwwL0SNDd6mCSEYJXw%2B9gDhBuJdCeHictbA3HmkcoiFWnxRGsESdkU16M6Vmn77knTs5qeVofZiaVjFsKV3yOQWcIS%2FK%2BUvBtrDDmCk1OCRgXGhvrtT62wPPI14jm%2FhmWSk5GUv6T0ZRPB8vqzBZlocGPAHamgxIw7oRPjpaiZ8d9w3Myr4UKD181LoMDF8CVg4EYzlS%2FHYSkBXIpjj5oLJiBo44tyPufRHFbwhzFbV%2BD18bNJLSyoC1VhNnI4IPd631VjCQerYOmOd5XO%2BBRQaUrbueM9pEtcnOkLaHnMjSs5bpEnwrAasGxNTXNy194kMt9TRdSCRHy%2FxX9HxgjNBSsrqH%2FcpGXNcD5KcdqAv4VYdoD%2FN%2FNJmAOXHB45hChJseROOgsxS59gJtFBKwEgP53AkdmUCIwj4jGI5zPV4wBw6vDEp4Tzq9cZ8S3UVarKA7E7eZwlv42plRUqBkLBLoJlamXlb1IrHB8MUputIwt9GlvwY6pQHti%2Bp%2F5lkx%2FknKSxaguyBNWhzhkYXtbfHkO4JflVjx7mnrVBOFxZnEB4xT5%2Bt%2FkeNPhONL%2BdxmgoROEwDDGyw%2FyC1Tbqo%2Fu%2FJy1qS1VgoP8dd5g8L71ClJugxPabxV1aBNhXxamJEH5PNkBBwLmZh%2BawcFaQdifFc7k3qzPe5Fz89GimhOBRr3YDxNvSjkmxnu2%2BjMKEuFDWIKSqeLuQRoIdq%2BGhg%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2004%3A00%3A14%20GMT&X-Amz-Signature=7dc3b1a298770b74a089eb23ad3e4deeff3b31d8b951d00bdd718b8b00fcc45f'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHP5YWAN5%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T200014Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCu5FJXG3S5aKeF7FyN2EBzNWsS9JBYsFDsYKJGuEbZEQIgAhZCaqLzvYqjyzdSgwHh1ncWapjmeSRzcn9wOTShDXgq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDPKRBg3%2B7gik4%2FdoJCq6A60R3mr2sOyO17YUbq%2Bv4aDc365v9rrCE13OkzA4%2BfkMac6Rf4uuQROrUAtfb%2Bjp%2BN1bGN7HKqLD9xOmXEdNljVWjute3j05HX1X%2BwwL0SNDd6mCSEYJXw%2B9gDhBuJdCeHictbA3HmkcoiFWnxRGsESdkU16M6Vmn77knTs5qeVofZiaVjFsKV3yOQWcIS%2FK%2BUvBtrDDmCk1OCRgXGhvrtT62wPPI14jm%2FhmWSk5GUv6T0ZRPB8vqzBZlocGPAHamgxIw7oRPjpaiZ8d9w3Myr4UKD181LoMDF8CVg4EYzlS%2FHYSkBXIpjj5oLJiBo44tyPufRHFbwhzFbV%2BD18bNJLSyoC1VhNnI4IPd631VjCQerYOmOd5XO%2BBRQaUrbueM9pEtcnOkLaHnMjSs5bpEnwrAasGxNTXNy194kMt9TRdSCRHy%2FxX9HxgjNBSsrqH%2FcpGXNcD5KcdqAv4VYdoD%2FN%2FNJmAOXHB45hChJseROOgsxS59gJtFBKwEgP53AkdmUCIwj4jGI5zPV4wBw6vDEp4Tzq9cZ8S3UVarKA7E7eZwlv42plRUqBkLBLoJlamXlb1IrHB8MUputIwt9GlvwY6pQHti%2Bp%2F5lkx%2FknKSxaguyBNWhzhkYXtbfHkO4JflVjx7mnrVBOFxZnEB4xT5%2Bt%2FkeNPhONL%2BdxmgoROEwDDGyw%2FyC1Tbqo%2Fu%2FJy1qS1VgoP8dd5g8L71ClJugxPabxV1aBNhXxamJEH5PNkBBwLmZh%2BawcFaQdifFc7k3qzPe5Fz89GimhOBRr3YDxNvSjkmxnu2%2BjMKEuFDWIKSqeLuQRoIdq%2BGhg%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2004%3A00%3A14%20GMT&X-Amz-Signature=0f9ce99e56cb03f770c50e9cb077729c2984aca2ddc391f3951e30517699c9fc'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHP5YWAN5%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T200014Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCu5FJXG3S5aKeF7FyN2EBzNWsS9JBYsFDsYKJGuEbZEQIgAhZCaqLzvYqjyzdSgwHh1ncWapjmeSRzcn9wOTShDXgq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDPKRBg3%2B7gik4%2FdoJCq6A60R3mr2sOyO17YUbq%2Bv4aDc365v9rrCE13OkzA4%2BfkMac6Rf4uuQROrUAtfb%2Bjp%2BN1bGN7HKqLD9xOmXEdNljVWjute3j05HX1X%2BwwL0SNDd6mCSEYJXw%2B9gDhBuJdCeHictbA3HmkcoiFWnxRGsESdkU16M6Vmn77knTs5qeVofZiaVjFsKV3yOQWcIS%2FK%2BUvBtrDDmCk1OCRgXGhvrtT62wPPI14jm%2FhmWSk5GUv6T0ZRPB8vqzBZlocGPAHamgxIw7oRPjpaiZ8d9w3Myr4UKD181LoMDF8CVg4EYzlS%2FHYSkBXIpjj5oLJiBo44tyPufRHFbwhzFbV%2BD18bNJLSyoC1VhNnI4IPd631VjCQerYOmOd5XO%2BBRQaUrbueM9pEtcnOkLaHnMjSs5bpEnwrAasGxNTXNy194kMt9TRdSCRHy%2FxX9HxgjNBSsrqH%2FcpGXNcD5KcdqAv4VYdoD%2FN%2FNJmAOXHB45hChJseROOgsxS59gJtFBKwEgP53AkdmUCIwj4jGI5zPV4wBw6vDEp4Tzq9cZ8S3UVarKA7E7eZwlv42plRUqBkLBLoJlamXlb1IrHB8MUputIwt9GlvwY6pQHti%2Bp%2F5lkx%2FknKSxaguyBNWhzhkYXtbfHkO4JflVjx7mnrVBOFxZnEB4xT5%2Bt%2FkeNPhONL%2BdxmgoROEwDDGyw%2FyC1Tbqo%2Fu%2FJy1qS1VgoP8dd5g8L71ClJugxPabxV1aBNhXxamJEH5PNkBBwLmZh%2BawcFaQdifFc7k3qzPe5Fz89GimhOBRr3YDxNvSjkmxnu2%2BjMKEuFDWIKSqeLuQRoIdq%2BGhg%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2004%3A00%3A14%20GMT&X-Amz-Signature=da1c4895f2fa3a2f2fe4f53bcf4a68e532af0f496f318a44eac9fac4b6318805'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHP5YWAN5%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T200014Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCu5FJXG3S5aKeF7FyN2EBzNWsS9JBYsFDsYKJGuEbZEQIgAhZCaqLzvYqjyzdSgwHh1ncWapjmeSRzcn9wOTShDXgq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDPKRBg3%2B7gik4%2FdoJCq6A60R3mr2sOyO17YUbq%2Bv4aDc365v9rrCE13OkzA4%2BfkMac6Rf4uuQROrUAtfb%2Bjp%2BN1bGN7HKqLD9xOmXEdNljVWjute3j05HX1X%2BwwL0SNDd6mCSEYJXw%2B9gDhBuJdCeHictbA3HmkcoiFWnxRGsESdkU16M6Vmn77knTs5qeVofZiaVjFsKV3yOQWcIS%2FK%2BUvBtrDDmCk1OCRgXGhvrtT62wPPI14jm%2FhmWSk5GUv6T0ZRPB8vqzBZlocGPAHamgxIw7oRPjpaiZ8d9w3Myr4UKD181LoMDF8CVg4EYzlS%2FHYSkBXIpjj5oLJiBo44tyPufRHFbwhzFbV%2BD18bNJLSyoC1VhNnI4IPd631VjCQerYOmOd5XO%2BBRQaUrbueM9pEtcnOkLaHnMjSs5bpEnwrAasGxNTXNy194kMt9TRdSCRHy%2FxX9HxgjNBSsrqH%2FcpGXNcD5KcdqAv4VYdoD%2FN%2FNJmAOXHB45hChJseROOgsxS59gJtFBKwEgP53AkdmUCIwj4jGI5zPV4wBw6vDEp4Tzq9cZ8S3UVarKA7E7eZwlv42plRUqBkLBLoJlamXlb1IrHB8MUputIwt9GlvwY6pQHti%2Bp%2F5lkx%2FknKSxaguyBNWhzhkYXtbfHkO4JflVjx7mnrVBOFxZnEB4xT5%2Bt%2FkeNPhONL%2BdxmgoROEwDDGyw%2FyC1Tbqo%2Fu%2FJy1qS1VgoP8dd5g8L71ClJugxPabxV1aBNhXxamJEH5PNkBBwLmZh%2BawcFaQdifFc7k3qzPe5Fz89GimhOBRr3YDxNvSjkmxnu2%2BjMKEuFDWIKSqeLuQRoIdq%2BGhg%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2004%3A00%3A14%20GMT&X-Amz-Signature=f5021e37b36103b9035d76dfc188c04f16b6adf8a5ca852695b4e4e13d1dfdd8'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHP5YWAN5%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T200014Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCu5FJXG3S5aKeF7FyN2EBzNWsS9JBYsFDsYKJGuEbZEQIgAhZCaqLzvYqjyzdSgwHh1ncWapjmeSRzcn9wOTShDXgq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDPKRBg3%2B7gik4%2FdoJCq6A60R3mr2sOyO17YUbq%2Bv4aDc365v9rrCE13OkzA4%2BfkMac6Rf4uuQROrUAtfb%2Bjp%2BN1bGN7HKqLD9xOmXEdNljVWjute3j05HX1X%2BwwL0SNDd6mCSEYJXw%2B9gDhBuJdCeHictbA3HmkcoiFWnxRGsESdkU16M6Vmn77knTs5qeVofZiaVjFsKV3yOQWcIS%2FK%2BUvBtrDDmCk1OCRgXGhvrtT62wPPI14jm%2FhmWSk5GUv6T0ZRPB8vqzBZlocGPAHamgxIw7oRPjpaiZ8d9w3Myr4UKD181LoMDF8CVg4EYzlS%2FHYSkBXIpjj5oLJiBo44tyPufRHFbwhzFbV%2BD18bNJLSyoC1VhNnI4IPd631VjCQerYOmOd5XO%2BBRQaUrbueM9pEtcnOkLaHnMjSs5bpEnwrAasGxNTXNy194kMt9TRdSCRHy%2FxX9HxgjNBSsrqH%2FcpGXNcD5KcdqAv4VYdoD%2FN%2FNJmAOXHB45hChJseROOgsxS59gJtFBKwEgP53AkdmUCIwj4jGI5zPV4wBw6vDEp4Tzq9cZ8S3UVarKA7E7eZwlv42plRUqBkLBLoJlamXlb1IrHB8MUputIwt9GlvwY6pQHti%2Bp%2F5lkx%2FknKSxaguyBNWhzhkYXtbfHkO4JflVjx7mnrVBOFxZnEB4xT5%2Bt%2FkeNPhONL%2BdxmgoROEwDDGyw%2FyC1Tbqo%2Fu%2FJy1qS1VgoP8dd5g8L71ClJugxPabxV1aBNhXxamJEH5PNkBBwLmZh%2BawcFaQdifFc7k3qzPe5Fz89GimhOBRr3YDxNvSjkmxnu2%2BjMKEuFDWIKSqeLuQRoIdq%2BGhg%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2004%3A00%3A14%20GMT&X-Amz-Signature=ecd123dadfe5bd958b01bb7d3fd4c1c662c0565da69aef11cec46fdffbca2527'
